$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Fitness (column C) for all data rows (2 through 252) to the constant value 7293
$ws.Range("C2:C252").Value = 7293
